$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert every cell in the table to Text format first, so that the new
#     values typed below are stored as literal text instead of being
#     auto-coerced back into numbers/dates. Columns G:H only have data in the
#     first three rows, so keep the touched range the same shape as before. ---
$ws.Range("A1:F17").NumberFormat = "@"
$ws.Range("G1:H3").NumberFormat = "@"

# --- Update the two fields whose content actually changed ---
# identidade (col E): strip the thousands separators so the id reads as a
# plain digit string (stored as text).
$ws.Range("E2").Value = "1122233334"
$ws.Range("E3").Value = "4433322221"

# admissao (col F): replace the date values with plain text dates.
$ws.Range("F2").Value = "21/02/2020"
$ws.Range("F3").Value = "15/04/2020"

# --- Selection moves from H3 to F3 ---
$ws.Range("F3").Select()
